# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E17:E37) is re-sorted into ascending order
# (it previously ran descending from 2206 down to 2010; it now runs
# ascending from 2010 up to 2206). As a side effect of the underlying
# data refresh, the "Valor Mora" amounts in F17 and F37 swap places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2010","2011","2012","2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112","2201","2202","2203","2204","2205","2206")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 17 + $i
    $ws.Cells.Item($row, 5).Value2 = $periods[$i]
}

# Valor Mora amounts for the first/last rows swap.
$ws.Cells.Item(17, 6).Value2 = 19897
$ws.Cells.Item(37, 6).Value2 = 25749
